$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D (bold, matching B1/C1 style)
$ws.Range("D1").Value = "Kalman with Eigen C++ (AvgTime for 1000 Runs)"
$ws.Range("D1").Font.Bold = $true

# Match column width of column C (closest achievable quantized width)
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# Fill in the new data values for rows 2-11
$values = @(48.85, 57.07, 61.999, 59.22, 62.827, 56.926, 58.397, 55.204, 56.714, 56.73)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $values[$i]
}

# Add the AVERAGE formula in D12
$ws.Range("D12").Formula = "=AVERAGE(D2:D11)"
